$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 previously only held a leftover test case name ("verifyListBox") with
# no other data. Replace it with a full data row so the test cases can run
# in parallel (per commit message), and add a new row 3 entry.
$ws.Range("A3").Value = "launchOrangeHRM"
$ws.Range("B3").Value = "Admin"
$ws.Range("C3").Value = "admin123"
$ws.Range("D3").Value = "john"

# Update the active selection to reflect the newly added rows.
[void]$ws.Range("D4:D5").Select()
